# The document currently ends with the paragraph:
#   "Cuando vemos un inventario debe tener la opción para regresar a la lista de inventarios"
# We need to append three new list paragraphs (same "Prrafodelista" / numId 1 list style)
# after it:
#   1) "Si agregó más productos el inventario ya no se toca, ya es caso cerrado."
#   2) "Cuando sale una mercaderia por ejemplo por falla del producto se tiene que hacer
#       una nota de crédito por devolución del producto"  (with spell-check markers
#       wrapping the misspelled word "mercaderia")
#   3) an empty list paragraph (no run at all)

$d = $word.ActiveDocument

# Create a fresh paragraph at the end of the document to use as an anchor/insertion
# point. Word's COM layer seeds it with the same paragraph formatting (style, numPr,
# run-properties) as the paragraph it follows, which is exactly the formatting our new
# paragraphs need.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last

# Build the raw WordprocessingML for the two new text-bearing paragraphs. Inserting this
# OOXML fragment over the freshly-created (empty) anchor paragraph's Range replaces just
# that paragraph, splicing in both new paragraphs while Word automatically re-creates a
# trailing paragraph mark (carrying the same pPr/rPr, with no run) to preserve the
# document's end-of-content position -- i.e. exactly the empty third paragraph the diff
# calls for.
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$p1 = "<w:p><w:pPr><w:pStyle w:val='Prrafodelista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='es-MX'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='es-MX'/></w:rPr><w:t>Si agregó más productos el inventario ya no se toca, ya es caso cerrado.</w:t></w:r></w:p>"

$p2 = "<w:p><w:pPr><w:pStyle w:val='Prrafodelista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='es-MX'/></w:rPr></w:pPr>" + `
        "<w:r><w:rPr><w:lang w:val='es-MX'/></w:rPr><w:t xml:space='preserve'>Cuando sale una </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r><w:rPr><w:lang w:val='es-MX'/></w:rPr><w:t>mercaderia</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r><w:rPr><w:lang w:val='es-MX'/></w:rPr><w:t xml:space='preserve'> por ejemplo por falla del producto se tiene que hacer una nota de crédito por devolución del producto</w:t></w:r>" + `
        "</w:p>"

$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" + `
         "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" + `
           "<pkg:xmlData><w:document $ns><w:body>$p1$p2</w:body></w:document></pkg:xmlData>" + `
         "</pkg:part></pkg:package>"

[void]$anchor.Range.InsertXML($xml)
